$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused translation column (C)
$ws.Columns.Item(3).Delete()

# Fix typo in the "curva esquerda agressiva" label (row 2, column A)
$ws.Cells.Item(2, 1).Value = "curva_esquerda_agressiva"

# Insert a new header row at the top and shift existing data down
$ws.Rows.Item(1).Insert()
$ws.Cells.Item(1, 1).Value = "event"
$ws.Cells.Item(1, 2).Value = "val"

# Add a new data row at the bottom
$ws.Cells.Item(8, 1).Value = "troca_faixa_esquerda_agressiva"
$ws.Cells.Item(8, 2).Value = 4

[void]$ws.Range("B8").Select()
